# Auto-generated edit script applying cached-value corrections
# from the Coeurl_Profits market-data refresh (scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2441688.8
$ws.Range("J17").Value = 2790302.5
$ws.Range("L17").Value = 8370907.5
$ws.Range("N17").Value = -8371243.5
$ws.Range("H41").Value = 1071.2413
$ws.Range("J41").Value = 493.25
$ws.Range("L41").Value = 493.25
$ws.Range("N41").Value = -1373.25
$ws.Range("H70").Value = 3778.0476
$ws.Range("I70").Value = 4564.5386
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 13693.6158
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -13423.6158
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 3778.0476
$ws.Range("I73").Value = 4564.5386
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 13693.6158
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -12757.6158
$ws.Range("N73").Value = -9372
$ws.Range("H86").Value = 3716.9443
$ws.Range("I86").Value = 2873.7
$ws.Range("J86").Value = 4771
$ws.Range("K86").Value = 2873.7
$ws.Range("L86").Value = 4771
$ws.Range("M86").Value = -1750.7
$ws.Range("N86").Value = -7017
$ws.Range("H89").Value = 3716.9443
$ws.Range("I89").Value = 2873.7
$ws.Range("J89").Value = 4771
$ws.Range("K89").Value = 14368.5
$ws.Range("L89").Value = 23855
$ws.Range("M89").Value = -8752.5
$ws.Range("N89").Value = -35087
$ws.Range("H106").Value = 2759.6667
$ws.Range("I106").Value = 2479.625
$ws.Range("K106").Value = 2479.625
$ws.Range("M106").Value = -1848.625
$ws.Range("H112").Value = 40821.855
$ws.Range("I112").Value = 2100
$ws.Range("K112").Value = 6300
$ws.Range("M112").Value = -5192

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5870.75
$ws.Range("I32").Value = 2921.9658
$ws.Range("K32").Value = 2921.9658
$ws.Range("M32").Value = -2634.9658
$ws.Range("H43").Value = 22425
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("N43").Value = -25626
$ws.Range("H45").Value = 9491.385
$ws.Range("J45").Value = 2019.8
$ws.Range("L45").Value = 2019.8
$ws.Range("N45").Value = -2773.8
$ws.Range("H101").Value = 39985
$ws.Range("J101").Value = 39985
$ws.Range("L101").Value = 39985
$ws.Range("N101").Value = -46475
$ws.Range("H104").Value = 121797.1
$ws.Range("J104").Value = 121797.1
$ws.Range("L104").Value = 121797.1
$ws.Range("N104").Value = -128785.1
$ws.Range("H106").Value = 19666.334
$ws.Range("J106").Value = 19666.334
$ws.Range("L106").Value = 19666.334
$ws.Range("N106").Value = -22190.334
$ws.Range("H110").Value = 7912.8
$ws.Range("I110").Value = 8646.929
$ws.Range("K110").Value = 8646.929
$ws.Range("M110").Value = -6601.929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 93340.73
$ws.Range("I31").Value = 126972.375
$ws.Range("J31").Value = 3656.3333
$ws.Range("K31").Value = 126972.375
$ws.Range("L31").Value = 3656.3333
$ws.Range("M31").Value = -126677.375
$ws.Range("N31").Value = -4246.3333
$ws.Range("H34").Value = 93340.73
$ws.Range("I34").Value = 126972.375
$ws.Range("J34").Value = 3656.3333
$ws.Range("K34").Value = 126972.375
$ws.Range("L34").Value = 3656.3333
$ws.Range("M34").Value = -126770.375
$ws.Range("N34").Value = -4060.3333
$ws.Range("H86").Value = 8652.968000000001
$ws.Range("I86").Value = 9676.083000000001
$ws.Range("K86").Value = 9676.083000000001
$ws.Range("M86").Value = -8553.083000000001
$ws.Range("H88").Value = 16566.5
$ws.Range("J88").Value = 16566.5
$ws.Range("L88").Value = 16566.5
$ws.Range("N88").Value = -17378.5
$ws.Range("H89").Value = 8652.968000000001
$ws.Range("I89").Value = 9676.083000000001
$ws.Range("K89").Value = 48380.415
$ws.Range("M89").Value = -42764.415
$ws.Range("H91").Value = 16566.5
$ws.Range("J91").Value = 16566.5
$ws.Range("L91").Value = 16566.5
$ws.Range("N91").Value = -19374.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1615.5769
$ws.Range("I26").Value = 1305.3334
$ws.Range("J26").Value = 1779.8235
$ws.Range("K26").Value = 3916.0002
$ws.Range("L26").Value = 5339.470499999999
$ws.Range("M26").Value = -3628.0002
$ws.Range("N26").Value = -5915.470499999999
$ws.Range("H121").Value = 818.7273
$ws.Range("I121").Value = 282
$ws.Range("J121").Value = 1266
$ws.Range("K121").Value = 846
$ws.Range("L121").Value = 3798
$ws.Range("M121").Value = 464
$ws.Range("N121").Value = -6418
$ws.Range("H122").Value = 1281.95
$ws.Range("J122").Value = 1793.4286
$ws.Range("L122").Value = 16140.8574
$ws.Range("N122").Value = -21040.8574
$ws.Range("H131").Value = 20156.793
$ws.Range("I131").Value = 77875.234
$ws.Range("K131").Value = 233625.702
$ws.Range("M131").Value = -228585.702
$ws.Range("H132").Value = 1593.05
$ws.Range("I132").Value = 1295.2963
$ws.Range("J132").Value = 2211.4614
$ws.Range("K132").Value = 11657.6667
$ws.Range("L132").Value = 19903.1526
$ws.Range("M132").Value = -9127.6667
$ws.Range("N132").Value = -24963.1526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 42860000
$ws.Range("I10").Value = 53500000
$ws.Range("J10").Value = 299999
$ws.Range("K10").Value = 53500000
$ws.Range("L10").Value = 299999
$ws.Range("M10").Value = -53499831
$ws.Range("N10").Value = -300337
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 14800.0625
$ws.Range("I70").Value = 12531.909
$ws.Range("K70").Value = 12531.909
$ws.Range("M70").Value = -12261.909
$ws.Range("H73").Value = 14800.0625
$ws.Range("I73").Value = 12531.909
$ws.Range("K73").Value = 12531.909
$ws.Range("M73").Value = -11595.909
$ws.Range("H107").Value = 437.19354
$ws.Range("J107").Value = 309.16666
$ws.Range("L107").Value = 309.16666
$ws.Range("N107").Value = -4149.16666
$ws.Range("H126").Value = 2732.6191
$ws.Range("I126").Value = 2667.0625
$ws.Range("J126").Value = 2942.4
$ws.Range("K126").Value = 8001.1875
$ws.Range("L126").Value = 8827.200000000001
$ws.Range("M126").Value = -5531.1875
$ws.Range("N126").Value = -13767.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5275.4443
$ws.Range("I7").Value = 3168.6428
$ws.Range("K7").Value = 3168.6428
$ws.Range("M7").Value = -3056.6428
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 1000
$ws.Range("K11").Value = 1000
$ws.Range("M11").Value = -860
$ws.Range("H24").Value = 13333
$ws.Range("J24").Value = 13333
$ws.Range("L24").Value = 13333
$ws.Range("N24").Value = -14019
$ws.Range("H61").Value = 2473.0833
$ws.Range("I61").Value = 2276.261
$ws.Range("K61").Value = 2276.261
$ws.Range("M61").Value = -2074.261
$ws.Range("H100").Value = 223967.11
$ws.Range("I100").Value = 1375.5
$ws.Range("J100").Value = 402040.4
$ws.Range("K100").Value = 1375.5
$ws.Range("L100").Value = 402040.4
$ws.Range("M100").Value = -834.5
$ws.Range("N100").Value = -403122.4
$ws.Range("H113").Value = 2473.0833
$ws.Range("I113").Value = 2276.261
$ws.Range("K113").Value = 2276.261
$ws.Range("M113").Value = -106.261
$ws.Range("H122").Value = 5113.276
$ws.Range("I122").Value = 4935.4116
$ws.Range("K122").Value = 14806.2348
$ws.Range("M122").Value = -12356.2348
$ws.Range("H126").Value = 5275.4443
$ws.Range("I126").Value = 3168.6428
$ws.Range("K126").Value = 9505.928400000001
$ws.Range("M126").Value = -7035.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2354.6
$ws.Range("I4").Value = 681.25
$ws.Range("J4").Value = 3470.1667
$ws.Range("K4").Value = 681.25
$ws.Range("L4").Value = 3470.1667
$ws.Range("M4").Value = -568.25
$ws.Range("N4").Value = -3696.1667
$ws.Range("H9").Value = 26803.2
$ws.Range("J9").Value = 15004.5
$ws.Range("L9").Value = 15004.5
$ws.Range("N9").Value = -15284.5
$ws.Range("H122").Value = 2278.4358
$ws.Range("I122").Value = 2116.2856
$ws.Range("K122").Value = 6348.8568
$ws.Range("M122").Value = -3898.8568
$ws.Range("H136").Value = 1954.6842
$ws.Range("I136").Value = 1832.3529
$ws.Range("J136").Value = 2994.5
$ws.Range("K136").Value = 5497.0587
$ws.Range("L136").Value = 8983.5
$ws.Range("M136").Value = -2947.0587
$ws.Range("N136").Value = -14083.5

